$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Compatibility" column (filter results now show supported platforms)
$ws.Range("E1").Value = "Compatibility"
$ws.Range("E2").Value = "Windows"
$ws.Range("E3").Value = "macOS"
$ws.Range("E4").Value = "SteamOS + Linux"

# David's entry now sorts by lowest price instead of relevance
$ws.Range("D4").Value = "Lowest Price"

$ws.Range("C5").Select()
